# Updates D (Price) and E (Volume/1h change) columns for the cryptos
# sheet, row by row, matching the scraped-data refresh.
#
# Price cells (column D) are stored as literal text in the workbook (they
# use '.'-as-thousands-separator formatting, e.g. "37.439.29"). Several of
# the new prices parse as plain numbers (e.g. "0.998"), so a leading "'"
# is used on those to force Excel to keep them as text instead of silently
# re-typing the cell as a number (exactly like typing into the Excel UI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.439.29"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "2.050.08"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "'228.49"
$ws.Range("E5").Value = "  -2.22%  "

$ws.Range("E6").Value = "  -1.95%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'56.16"
$ws.Range("E8").Value = "  -3.42%  "

$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("D10").Value = "'0.0804"
$ws.Range("E10").Value = "  +2.89%  "

$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").Value = "2.352.90"
$ws.Range("E12").Value = "  -1.75%  "

$ws.Range("D13").Value = "'14.51"
$ws.Range("E13").Value = "  -4.93%  "

$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("E15").Value = "  -3.08%  "

$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("D17").Value = "2.052.10"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "37.326.71"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("D19").Value = "'6.07"
$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").Value = "'69.76"
$ws.Range("E20").Value = "  -1.68%  "

$ws.Range("D21").Value = "0.0₃0846"
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("D22").Value = "'225.63"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "  -4.54%  "

$ws.Range("D26").Value = "'9.48"
$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("D27").Value = "'168.25"
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("E28").Value = "  -4.70%  "

$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").Value = "'18.89"
$ws.Range("E30").Value = "  -3.07%  "

$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("E32").Value = "  -3.32%  "

$ws.Range("D33").Value = "'0.0612"
$ws.Range("E33").Value = "  -3.11%  "

$ws.Range("D34").Value = "'4.52"
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("E35").Value = "  -3.91%  "

$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "'3.19"
$ws.Range("E38").Value = "  -3.98%  "

$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").Value = "'0.0220"
$ws.Range("E40").Value = "  -5.93%  "

$ws.Range("D41").Value = "1.496.54"
$ws.Range("E41").Value = "  +3.22%  "

$ws.Range("D42").Value = "'2.86"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").Value = "'16.76"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").Value = "'95.96"
$ws.Range("E44").Value = "  -5.11%  "

$ws.Range("D45").Value = "'0.0932"
$ws.Range("E45").Value = "  -3.78%  "

$ws.Range("E46").Value = "  -4.11%  "

$ws.Range("D47").Value = "'1.01"
$ws.Range("E47").Value = "  -4.19%  "

$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "  -1.36%  "

$ws.Range("E50").Value = "  -9.07%  "

$ws.Range("D51").Value = "2.237.46"
$ws.Range("E51").Value = "  -1.84%  "
